# Update the "Indexes" sheet: rename step labels and fix the numeric values
# to match the new naming convention (start=0, stop=-1, step=1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indexes")

$ws.Range("B2").Value = "step_1_start"
$ws.Range("B3").Value = "step_1_stop"
$ws.Range("B4").Value = "step_1_step"
$ws.Range("B5").Value = "step_-1_start"
$ws.Range("B6").Value = "step_-1_stop"
$ws.Range("B7").Value = "step_-1_step"

$ws.Range("C2").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("C7").Value = 1
